# Applies the crypto-price / volume(1h) refresh captured by the commit's
# XML diff. Columns B/C/D/E hold Coin / Link / Price / Volume(1h); D's
# 'prices' are free-form text (e.g. '29.317.07', '0.00001020') rather than
# real numbers, so a plain .Value assignment would let Excel re-interpret
# them numerically and mangle the literal text (dropped trailing zeros,
# float noise, '.' used as a thousands separator, etc). Forcing the cell to
# Text format ('@') before the assignment keeps the exact string, then
# ClearFormats() drops that temporary formatting again so the cell's style
# ends up identical to how it started (no lingering NumberFormat change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.317.07'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.14%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.839.44'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.30%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9985'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.16'
$ws.Range('D5').ClearFormats()

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6272'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.09%  '

# Row 7
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07427'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.92%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2892'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.32%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.93'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.27%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07717'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.02%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.844.45'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.04%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.955'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.83%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6741'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.75%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001020'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.98%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.69'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.49%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.221'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.14%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.281.42'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.40%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.82'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.29%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.28'
$ws.Range('D20').ClearFormats()

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9997'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.09%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.342'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.62%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.19%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '157.97'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.46%  '

# Row 25
$ws.Range('E25').Value = '  +0.76%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1347'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.09%  '

# Row 27
$ws.Range('E27').Value = '  -1.09%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.07261'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +12.98%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.457'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.12%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.478'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.42%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.038'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.34%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.040'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.59%  '

# Row 33
$ws.Range('E33').Value = '  -0.49%  '

# Row 34
$ws.Range('E34').Value = '  -0.03%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6945'
$ws.Range('D35').ClearFormats()

# Row 36
$ws.Range('E36').Value = '  -0.19%  '

# Row 37
$ws.Range('E37').Value = '  +0.54%  '

# Row 38
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.812'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.82%  '

# Row 39
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.871'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.12%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.234.83'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.94%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9346'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.70%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9997'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.11%  '

# Row 43
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.984.61'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.08%  '

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.99'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.40%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.31'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.19%  '

# Row 46
$ws.Range('E46').Value = '  +1.76%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.702'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.38%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.934'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.05%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1138'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.11%  '

# Row 50
$ws.Range('E50').Value = '  -1.65%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3899'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.12%  '
